$wb = $excel.ActiveWorkbook

$credits = $wb.Worksheets.Item("Credits")
$assets = $wb.Worksheets.Item("Assets")

# Add new row to Credits sheet
$credits.Range("A26").Value = "Environment"
$credits.Range("E26").Value = "https://www.youtube.com/watch?v=Dk_0rf2YYCw"
$credits.Range("B26").Value = "Outdoor Sound"
$credits.Range("C26").Value = "Audio"
$credits.Range("D26").Value = ".mp3"
$credits.Range("F26").Value = "Free / No Copyright"

# Update Assets sheet selection (while active) before switching away
[void]$assets.Activate()
[void]$assets.Range("A12").Select()

# Switch active sheet to Credits and update view state
[void]$credits.Activate()
$excel.ActiveWindow.Zoom = 85
[void]$credits.Range("E29").Select()
